# "Uncommented RAD Extension Payments Code and Test Data."
# The Extension Payments test row (row 3) had its Execute flag ("C3")
# disabled with the placeholder "DONOTRUN" value; re-enable it by setting
# it back to "Y", matching the Execute column of the other (enabled) rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C3").Value = "Y"

# Reflect the cell the author was working in when the sheet was saved.
$ws.Range("C3").Select()
